$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 45-46; this shifts the previous rows 45-88
# down to 47-90, matching the target dimension A1:R90.
$ws.Range("A45:A46").EntireRow.Insert()

# New row 45 data
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 44895
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 300000000
$ws.Range("G45").Value = "Espárragos"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = 1100
$ws.Range("N45").Value = '$/kilo'
$ws.Range("O45").Value = "Provincia de Linares"
$ws.Range("P45").Value = 1100
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"

# New row 46 data
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 44895
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = 300000000
$ws.Range("G46").Value = "Espárragos"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 1100
$ws.Range("L46").Value = 1100
$ws.Range("M46").Value = 1100
$ws.Range("N46").Value = '$/kilo'
$ws.Range("O46").Value = "Provincia de Linares"
$ws.Range("P46").Value = 1100
$ws.Range("Q46").Value = 1
$ws.Range("R46").Value = "Hortaliza"
